# Update the LR-pairs worksheet with the new TPM-derived data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (old rows 5, 6 and 7 are no longer part of the
# result set after recomputing with the new TPM values / renamed cluster).
$ws.Rows("5:7").Delete()

# Row 2: ECs -> FAPs (unchanged pairing, refreshed specificity/weight values)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Lgr6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06762866666666667
$ws.Range("H2").Value = 0.202886
$ws.Range("I2").Value = 0.0134153952845566
$ws.Range("J2").Value = 0.0134153952845566
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1753453333333333
$ws.Range("N2").Value = 0.5260359999999999
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.01185837109955556
$ws.Range("R2").Value = 0.106725339896
$ws.Range("S2").Value = 0.0134153952845566
$ws.Range("T2").Value = 0.0134153952845566

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.954393
$ws.Range("H3").Value = 14.863179
$ws.Range("I3").Value = 0.9827953701592058
$ws.Range("J3").Value = 0.9827953701592059
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1753453333333333
$ws.Range("N3").Value = 0.5260359999999999
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.8687296920493333
$ws.Range("R3").Value = 7.818567228443999
$ws.Range("S3").Value = 0.9827953701592058
$ws.Range("T3").Value = 0.9827953701592059

# Row 4: Resolving-Mac (formerly "MuSCs") -> FAPs
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lgr6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.019102
$ws.Range("H4").Value = 0.057306
$ws.Range("I4").Value = 0.003789234556237495
$ws.Range("J4").Value = 0.003789234556237496
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1753453333333333
$ws.Range("N4").Value = 0.5260359999999999
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.003349446557333333
$ws.Range("R4").Value = 0.030145019016
$ws.Range("S4").Value = 0.003789234556237495
$ws.Range("T4").Value = 0.003789234556237496
